# Update the "Films" demo worksheet: insert a "Genre" column between Name and
# Year, and expand the sample data to the full film list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 = headers, rows 2-10 = film data: Name | Genre | Year
$rows = @(
    @("Name",                         "Genre",     "Year"),
    @("Batman - The Dark Knight",     "Fantasy",   2008),
    @("Batman Begins",                "Fantasy",   2005),
    @("Cars",                         "Animation", 2006),
    @("Cars 2",                       "Animation", 2011),
    @("City of God",                  "Drama",     2003),
    @("Cool Runnings",                "Comedy",    1993),
    @("Fast and the Furious",         "Action",    2001),
    @("Iron Man",                     "Fantasy",   2008),
    @("Monty Python's Life of Brian", "Comedy",    1979)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $ws.Range("A" + ($r + 1)).Value = $rows[$r][0]
    $ws.Range("B" + ($r + 1)).Value = $rows[$r][1]
    $ws.Range("C" + ($r + 1)).Value = $rows[$r][2]
}

# Columns were best-fit/auto-sized to the new content.
$ws.Columns.Item(1).ColumnWidth = 25.45
$ws.Columns.Item(2).ColumnWidth = 6.8
$ws.Columns.Item(3).ColumnWidth = 4.15

# Selection ends up just past the last populated cell, like after typing the
# final entry and pressing Enter.
[void]$ws.Range("C11").Select()
